$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "checklist"
$ws.Name = "checklist"

# Insert a new row at position 7 (pushes the existing rows 7-8 down to 8-9),
# carrying the row formatting along (matches Excel's native row-insert behaviour).
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new "Ghost Hunter RPG02 Replay" entry.
$ws.Range("A7").Value = 2002
$ws.Range("B7").Value = "ゴーストハンターRPG02リプレイ 黒き死の仮面―草壁健一郎の事件簿"
$ws.Range("C7").Value = "Ghost Hunter RPG02 Replay Black Death Mask-Kenichiro Kusakabe's Casebook"
$ws.Range("D7").Value = "Fujimi Shobo"
$ws.Range("E7").Value = "black-death-mask-replay2.jpg"
$ws.Range("F7").Value = "replay"

# Widen column B to fit the new (longer) Japanese title text.
$ws.Columns.Item(2).ColumnWidth = 72

# Move the active selection to A8 (matches the saved view state in the diff).
[void]$ws.Range("A8").Select()
